$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$boolFormat = '"TRUE";"TRUE";"FALSE"'

function Set-UniqueColumn($headerRow, $dataRows) {
    # Header cell for the new "Unique" column
    $ws.Cells.Item($headerRow, 12).Value = "Unique"

    foreach ($r in $dataRows) {
        # Existing cells in columns B, C, K become bold
        $ws.Cells.Item($r, 2).Font.Bold = $true
        $ws.Cells.Item($r, 3).Font.Bold = $true
        $ws.Cells.Item($r, 11).Font.Bold = $true

        # New L column data cell: boolean-looking text formatted with a custom
        # TRUE/FALSE number format, left aligned
        $cell = $ws.Cells.Item($r, 12)
        $cell.Value = "'FALSE"
        $cell.NumberFormat = $boolFormat
        $cell.HorizontalAlignment = -4131
    }
}

# First "SAMPLE_TYPE" properties table (header row 12, data rows 13-15)
Set-UniqueColumn 12 @(13, 14, 15)

# Second "SAMPLE_TYPE" properties table (header row 20, data rows 21-23)
Set-UniqueColumn 20 @(21, 22, 23)

# Match the saved selection state from the target workbook
$ws.Range("M19").Select()
